$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all existing hyperlinks; they will be re-created below (pointing at the
# corrected mailto targets) once the cell values have been rewritten.
$ws.Hyperlinks.Delete()

# --- Username / password reshuffle on the existing rows -------------------
# Row 3 and 4 usernames swap places.
$ws.Range("A3").Value = "Staginguser_1"
$ws.Range("A4").Value = "Staginguser_3"

# Rows 8, 9 get new usernames (Staginguser_35 / _36) and row 8's password
# becomes the Paragon@2024 mailto address (was Password@123456).
$ws.Range("A8").Value = "Staginguser_35"
$ws.Range("B8").Value = "Paragon@2024"
$ws.Range("A9").Value = "Staginguser_36"

# --- New rows 10-12, duplicating row 8's data with fresh usernames --------
$ws.Range("A10").Value = "Staginguser_37"
$ws.Range("B10").Value = "Paragon@2024"
$ws.Range("C10").Value = "Automated_Campaign_7"
$ws.Range("D10").Value = "CCR312318"
$ws.Range("E10").Value = "VAT20"
$ws.Range("F10").Value = "PO123129"
$ws.Range("G10").Value = 123
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 435
$ws.Range("J10").Value = "Paragon CC (Dagenham)"
$ws.Range("K10").Value = "ERN1231237"
$ws.Range("L10").Value = 45
$ws.Range("M10").Value = 123
$ws.Range("N10").Value = 312

$ws.Range("A11").Value = "Staginguser_38"
$ws.Range("B11").Value = "Paragon@2024"
$ws.Range("C11").Value = "Automated_Campaign_7"
$ws.Range("D11").Value = "CCR312318"
$ws.Range("E11").Value = "VAT20"
$ws.Range("F11").Value = "PO123129"
$ws.Range("G11").Value = 123
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 435
$ws.Range("J11").Value = "Paragon CC (Dagenham)"
$ws.Range("K11").Value = "ERN1231237"
$ws.Range("L11").Value = 45
$ws.Range("M11").Value = 123
$ws.Range("N11").Value = 312

$ws.Range("A12").Value = "Staginguser_39"
$ws.Range("B12").Value = "Paragon@2024"
$ws.Range("C12").Value = "Automated_Campaign_7"
$ws.Range("D12").Value = "CCR312318"
$ws.Range("E12").Value = "VAT20"
$ws.Range("F12").Value = "PO123129"
$ws.Range("G12").Value = 123
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 435
$ws.Range("J12").Value = "Paragon CC (Dagenham)"
$ws.Range("K12").Value = "ERN1231237"
$ws.Range("L12").Value = 45
$ws.Range("M12").Value = 123
$ws.Range("N12").Value = 312

# --- Row 5, 6, 7 usernames, done last so the shared-string table allocates
# Staginguser_15/16/17 right after Staginguser_35-39 (matches source order).
$ws.Range("A5").Value = "Staginguser_15"
$ws.Range("B5").Value = "Paragon@2024"
$ws.Range("A6").Value = "Staginguser_16"
$ws.Range("A7").Value = "Staginguser_17"

# --- Recreate hyperlinks for column B, in the same order as the target ----
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Password@123456")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:Paragon@2024")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Paragon@2024")

# Match the saved selection from the edit.
$ws.Range("B14").Select() | Out-Null
